$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2349.7083
$ws.Range("I129").Value = 460
$ws.Range("J129").Value = 2521.5
$ws.Range("K129").Value = 1380
$ws.Range("L129").Value = 7564.5
$ws.Range("M129").Value = 3620
$ws.Range("N129").Value = -17564.5
$ws.Range("H132").Value = 8174.6313
$ws.Range("I132").Value = 7464.875
$ws.Range("K132").Value = 22394.625
$ws.Range("M132").Value = -19864.625
$ws.Range("H137").Value = 3603.75
$ws.Range("I137").Value = 4464.5
$ws.Range("J137").Value = 3316.8333
$ws.Range("K137").Value = 13393.5
$ws.Range("L137").Value = 9950.499899999999
$ws.Range("M137").Value = -10843.5
$ws.Range("N137").Value = -15050.4999
$ws.Range("H138").Value = 2110.0513
$ws.Range("I138").Value = 894.4375
$ws.Range("J138").Value = 2955.6956
$ws.Range("K138").Value = 2683.3125
$ws.Range("L138").Value = 8867.086800000001
$ws.Range("M138").Value = 2456.6875
$ws.Range("N138").Value = -19147.0868

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 580.6
$ws.Range("I4").Value = 580.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 580.6
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -464.6
$ws.Range("H74").Value = 18003134
$ws.Range("I74").Value = 22502768
$ws.Range("J74").Value = 4605.6
$ws.Range("K74").Value = 22502768
$ws.Range("L74").Value = 4605.6
$ws.Range("M74").Value = -22501894
$ws.Range("N74").Value = -6353.6
$ws.Range("H77").Value = 18003134
$ws.Range("I77").Value = 22502768
$ws.Range("J77").Value = 4605.6
$ws.Range("K77").Value = 112513840
$ws.Range("L77").Value = 23028
$ws.Range("M77").Value = -112509472
$ws.Range("N77").Value = -31764
$ws.Range("H88").Value = 111113570
$ws.Range("I88").Value = 1416
$ws.Range("J88").Value = 250003740
$ws.Range("K88").Value = 1416
$ws.Range("L88").Value = 250003740
$ws.Range("M88").Value = -1010
$ws.Range("N88").Value = -250004552
$ws.Range("H91").Value = 111113570
$ws.Range("I91").Value = 1416
$ws.Range("J91").Value = 250003740
$ws.Range("K91").Value = 1416
$ws.Range("L91").Value = 250003740
$ws.Range("M91").Value = -12
$ws.Range("N91").Value = -250006548
$ws.Range("H97").Value = 1541.1786
$ws.Range("I97").Value = 1865.5454
$ws.Range("J97").Value = 351.83334
$ws.Range("K97").Value = 1865.5454
$ws.Range("L97").Value = 351.83334
$ws.Range("M97").Value = -1369.5454
$ws.Range("N97").Value = -1343.83334
$ws.Range("H123").Value = 53085.8
$ws.Range("J123").Value = 53085.8
$ws.Range("L123").Value = 53085.8
$ws.Range("N123").Value = -62885.8
$ws.Range("H132").Value = 1873.7826
$ws.Range("I132").Value = 1662.4667
$ws.Range("J132").Value = 2270
$ws.Range("K132").Value = 4987.4001
$ws.Range("L132").Value = 6810
$ws.Range("M132").Value = -2457.4001
$ws.Range("N132").Value = -11870

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 322.27274
$ws.Range("I64").Value = 385
$ws.Range("J64").Value = 270
$ws.Range("K64").Value = 385
$ws.Range("L64").Value = 270
$ws.Range("M64").Value = -160
$ws.Range("N64").Value = -720
$ws.Range("H67").Value = 322.27274
$ws.Range("I67").Value = 385
$ws.Range("J67").Value = 270
$ws.Range("K67").Value = 385
$ws.Range("L67").Value = 270
$ws.Range("M67").Value = 395
$ws.Range("N67").Value = -1830
$ws.Range("H86").Value = 5884065.5
$ws.Range("I86").Value = 7144519.5
$ws.Range("K86").Value = 7144519.5
$ws.Range("M86").Value = -7143396.5
$ws.Range("H89").Value = 5884065.5
$ws.Range("I89").Value = 7144519.5
$ws.Range("K89").Value = 35722597.5
$ws.Range("M89").Value = -35716981.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 354.58823
$ws.Range("I22").Value = 301.0909
$ws.Range("J22").Value = 452.66666
$ws.Range("K22").Value = 301.0909
$ws.Range("L22").Value = 452.66666
$ws.Range("M22").Value = 48.90910000000002
$ws.Range("N22").Value = -1152.66666
$ws.Range("H99").Value = 3256794.2
$ws.Range("I99").Value = 8937799
$ws.Range("J99").Value = 10506
$ws.Range("K99").Value = 8937799
$ws.Range("L99").Value = 10506
$ws.Range("M99").Value = -8936301
$ws.Range("N99").Value = -13502
$ws.Range("H126").Value = 3256794.2
$ws.Range("I126").Value = 8937799
$ws.Range("J126").Value = 10506
$ws.Range("K126").Value = 26813397
$ws.Range("L126").Value = 31518
$ws.Range("M126").Value = -26810927
$ws.Range("N126").Value = -36458

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 47624490
$ws.Range("I7").Value = 200000290
$ws.Range("J7").Value = 7053.125
$ws.Range("K7").Value = 600000870
$ws.Range("L7").Value = 21159.375
$ws.Range("M7").Value = -600000758
$ws.Range("N7").Value = -21383.375

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 29.4
$ws.Range("I2").Value = 26.75
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 26.75
$ws.Range("L2").Value = 40
$ws.Range("M2").Value = 86.25
$ws.Range("N2").Value = -266
$ws.Range("H80").Value = 2873.182
$ws.Range("I80").Value = 2926.875
$ws.Range("J80").Value = 2856
$ws.Range("K80").Value = 2926.875
$ws.Range("L80").Value = 2856
$ws.Range("M80").Value = -1928.875
$ws.Range("N80").Value = -4852
$ws.Range("H83").Value = 2873.182
$ws.Range("I83").Value = 2926.875
$ws.Range("J83").Value = 2856
$ws.Range("K83").Value = 14634.375
$ws.Range("L83").Value = 14280
$ws.Range("M83").Value = -9642.375
$ws.Range("N83").Value = -24264
$ws.Range("H97").Value = 830.4545000000001
$ws.Range("I97").Value = 907.2222
$ws.Range("J97").Value = 485
$ws.Range("K97").Value = 907.2222
$ws.Range("L97").Value = 485
$ws.Range("M97").Value = -411.2222
$ws.Range("N97").Value = -1477
$ws.Range("H122").Value = 2678.6924
$ws.Range("I122").Value = 3225.875
$ws.Range("J122").Value = 1803.2
$ws.Range("K122").Value = 9677.625
$ws.Range("L122").Value = 5409.6
$ws.Range("M122").Value = -7227.625
$ws.Range("N122").Value = -10309.6

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7491.909
$ws.Range("I62").Value = 6513.75
$ws.Range("J62").Value = 8050.857
$ws.Range("K62").Value = 6513.75
$ws.Range("L62").Value = 8050.857
$ws.Range("M62").Value = -5889.75
$ws.Range("N62").Value = -9298.857
$ws.Range("H65").Value = 7491.909
$ws.Range("I65").Value = 6513.75
$ws.Range("J65").Value = 8050.857
$ws.Range("K65").Value = 32568.75
$ws.Range("L65").Value = 40254.285
$ws.Range("M65").Value = -29448.75
$ws.Range("N65").Value = -46494.285
$ws.Range("H123").Value = 30142.9
$ws.Range("J123").Value = 30142.9
$ws.Range("L123").Value = 30142.9
$ws.Range("N123").Value = -39942.9
$ws.Range("H132").Value = 2156.9736
$ws.Range("I132").Value = 1885.0769
$ws.Range("J132").Value = 2746.0833
$ws.Range("K132").Value = 5655.2307
$ws.Range("L132").Value = 8238.249899999999
$ws.Range("M132").Value = -3125.2307
$ws.Range("N132").Value = -13298.2499
$ws.Range("H136").Value = 1290.5518
$ws.Range("I136").Value = 958.4286
$ws.Range("K136").Value = 2875.2858
$ws.Range("M136").Value = -325.2857999999997
